$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.066.03"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.487.76"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.Value = "'594.81"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "
$c = $ws.Range("D6")
$c.Value = "'182.93"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$c = $ws.Range("D7")
$c.Value = "'0.616"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.13%  "
$c = $ws.Range("D8")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "3.485.48"
$ws.Range("E9").Value = "  -0.94%  "
$c = $ws.Range("D10")
$c.Value = "'0.142"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "4.087.00"
$ws.Range("E13").Value = "  -0.97%  "
$c = $ws.Range("D14")
$c.Value = "'32.34"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "68.043.67"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").Value = "3.487.96"
$ws.Range("E18").Value = "  -0.83%  "
$c = $ws.Range("D19")
$c.Value = "'6.22"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.64%  "
$c = $ws.Range("D20")
$c.Value = "'14.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.91%  "
$c = $ws.Range("D21")
$c.Value = "'395.97"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.40%  "
$c = $ws.Range("D22")
$c.Value = "'7.96"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.61%  "
$c = $ws.Range("D23")
$c.Value = "'5.85"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.50%  "
$c = $ws.Range("D24")
$c.Value = "'0.541"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  -0.27%  "
$c = $ws.Range("D26")
$c.Value = "'72.33"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.81%  "
$c = $ws.Range("D27")
$c.Value = "'0.0000123"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("E28").Value = "  +0.84%  "
$c = $ws.Range("D29")
$c.Value = "'0.178"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "
$c = $ws.Range("D30")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "
$c = $ws.Range("D31")
$c.Value = "'6.15"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").Value = "  -3.21%  "
$c = $ws.Range("D34")
$c.Value = "'23.64"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.50%  "
$c = $ws.Range("D35")
$c.Value = "'7.38"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -5.16%  "
$c = $ws.Range("D38")
$c.Value = "'161.92"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.93%  "
$c = $ws.Range("D39")
$c.Value = "'0.893"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.65%  "
$c = $ws.Range("D40")
$c.Value = "'2.84"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.82%  "
$ws.Range("E41").Value = "  -4.22%  "
$c = $ws.Range("D42")
$c.Value = "'6.80"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("E43").Value = "  -1.66%  "
$c = $ws.Range("D44")
$c.Value = "'26.24"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.99%  "
$c = $ws.Range("D45")
$c.Value = "'0.0723"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "2.752.93"
$ws.Range("E46").Value = "  -2.07%  "
$c = $ws.Range("D47")
$c.Value = "'26.35"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.35%  "
$c = $ws.Range("D48")
$c.Value = "'41.52"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.13%  "
$c = $ws.Range("D49")
$c.Value = "'0.0300"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.75%  "
$c = $ws.Range("D50")
$c.Value = "'330.51"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.03%  "
$ws.Range("E51").Value = "  -4.22%  "
